# Auto-generated edit script: updates IFRS financial data cells per commit
# "error solve ifrs list" - corrects misaligned/erroneous financial figures
# for rows 2-6 (2014-2018 actuals) and clears stale duplicate/erroneous
# forecast rows 7-9 (2019E-2021E) data cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update financial figures
$ws.Range("D2").Value = 3543
$ws.Range("E2").Value = 882
$ws.Range("F2").Value = 882
$ws.Range("G2").Value = 887
$ws.Range("H2").Value = 789
$ws.Range("I2").Value = 782
$ws.Range("J2").Value = 8
$ws.Range("K2").Value = 8372
$ws.Range("L2").Value = 1797
$ws.Range("M2").Value = 6576
$ws.Range("N2").Value = 6576
$ws.Range("P2").Value = 272
$ws.Range("Q2").Value = 310
$ws.Range("R2").Value = -248
$ws.Range("S2").Value = -53
$ws.Range("T2").Value = 221
$ws.Range("U2").Value = 89
$ws.Range("V2").Value = 933
$ws.Range("W2").Value = 24.9
$ws.Range("X2").Value = 22.28
$ws.Range("Y2").Value = 12.56
$ws.Range("Z2").Value = 9.72
$ws.Range("AA2").Value = 27.32
$ws.Range("AB2").Value = 2270.56
$ws.Range("AC2").Value = 1435
$ws.Range("AD2").Value = 6.06
$ws.Range("AE2").Value = 12110
$ws.Range("AF2").Value = 0.72
$ws.Range("AG2").Value = 55
$ws.Range("AH2").Value = 0.63
$ws.Range("AI2").Value = 3.84
$ws.Range("AJ2").Value = 50880170
$ws.Range("O2").ClearContents()

# Row 3: update financial figures
$ws.Range("D3").Value = 3061
$ws.Range("E3").Value = 858
$ws.Range("F3").Value = 858
$ws.Range("G3").Value = 852
$ws.Range("H3").Value = 778
$ws.Range("I3").Value = 778
$ws.Range("K3").Value = 8831
$ws.Range("L3").Value = 1582
$ws.Range("M3").Value = 7249
$ws.Range("N3").Value = 7249
$ws.Range("P3").Value = 272
$ws.Range("Q3").Value = 396
$ws.Range("R3").Value = -64
$ws.Range("S3").Value = -252
$ws.Range("T3").Value = 57
$ws.Range("U3").Value = 339
$ws.Range("V3").Value = 735
$ws.Range("W3").Value = 28.02
$ws.Range("X3").Value = 25.41
$ws.Range("Y3").Value = 11.25
$ws.Range("Z3").Value = 9.04
$ws.Range("AA3").Value = 21.82
$ws.Range("AB3").Value = 2508.35
$ws.Range("AC3").Value = 1428
$ws.Range("AD3").Value = 5.41
$ws.Range("AE3").Value = 13349
$ws.Range("AF3").Value = 0.58
$ws.Range("AG3").Value = 55
$ws.Range("AH3").Value = 0.71
$ws.Range("AI3").Value = 3.86
$ws.Range("AJ3").Value = 50880170
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()

# Row 4: update financial figures
$ws.Range("D4").Value = 3214
$ws.Range("E4").Value = 1091
$ws.Range("F4").Value = 1091
$ws.Range("G4").Value = 1107
$ws.Range("H4").Value = 994
$ws.Range("I4").Value = 994
$ws.Range("K4").Value = 9706
$ws.Range("L4").Value = 1498
$ws.Range("M4").Value = 8208
$ws.Range("N4").Value = 8208
$ws.Range("P4").Value = 272
$ws.Range("Q4").Value = 224
$ws.Range("R4").Value = -144
$ws.Range("S4").Value = -169
$ws.Range("T4").Value = 63
$ws.Range("U4").Value = 161
$ws.Range("V4").Value = 602
$ws.Range("W4").Value = 33.93
$ws.Range("X4").Value = 30.92
$ws.Range("Y4").Value = 12.86
$ws.Range("Z4").Value = 10.72
$ws.Range("AA4").Value = 18.25
$ws.Range("AB4").Value = 2862.85
$ws.Range("AC4").Value = 1824
$ws.Range("AD4").Value = 4.43
$ws.Range("AE4").Value = 15116
$ws.Range("AF4").Value = 0.53
$ws.Range("AG4").Value = 75
$ws.Range("AH4").Value = 0.93
$ws.Range("AI4").Value = 4.11
$ws.Range("AJ4").Value = 50880170
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()

# Row 5: update financial figures
$ws.Range("D5").Value = 3140
$ws.Range("E5").Value = 734
$ws.Range("F5").Value = 734
$ws.Range("G5").Value = 933
$ws.Range("H5").Value = 818
$ws.Range("I5").Value = 818
$ws.Range("K5").Value = 10855
$ws.Range("L5").Value = 1684
$ws.Range("M5").Value = 9171
$ws.Range("N5").Value = 9171
$ws.Range("P5").Value = 286
$ws.Range("Q5").Value = -141
$ws.Range("R5").Value = 333
$ws.Range("S5").Value = -172
$ws.Range("T5").Value = 88
$ws.Range("U5").Value = -229
$ws.Range("V5").Value = 569
$ws.Range("W5").Value = 23.38
$ws.Range("X5").Value = 26.06
$ws.Range("Y5").Value = 9.42
$ws.Range("Z5").Value = 7.96
$ws.Range("AA5").Value = 18.36
$ws.Range("AB5").Value = 3079.95
$ws.Range("AC5").Value = 1490
$ws.Range("AD5").Value = 5.09
$ws.Range("AE5").Value = 16099
$ws.Range("AF5").Value = 0.47
$ws.Range("AG5").Value = 75
$ws.Range("AH5").Value = 0.99
$ws.Range("AI5").Value = 5.24
$ws.Range("AJ5").Value = 53543977
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6: update financial figures
$ws.Range("D6").Value = 13600
$ws.Range("E6").Value = 1436
$ws.Range("F6").Value = 1436
$ws.Range("G6").Value = 876
$ws.Range("H6").Value = 885
$ws.Range("I6").Value = 622
$ws.Range("K6").Value = 38251
$ws.Range("L6").Value = 19423
$ws.Range("M6").Value = 18828
$ws.Range("N6").Value = 9592
$ws.Range("P6").Value = 286
$ws.Range("Q6").Value = 1356
$ws.Range("R6").Value = -1483
$ws.Range("S6").Value = 1901
$ws.Range("T6").Value = 2119
$ws.Range("U6").Value = -763
$ws.Range("V6").Value = 12590
$ws.Range("W6").Value = 10.56
$ws.Range("X6").Value = 6.51
$ws.Range("Y6").Value = 6.63
$ws.Range("Z6").Value = 3.61
$ws.Range("AA6").Value = 103.16
$ws.Range("AB6").Value = 3270.21
$ws.Range("AC6").Value = 1088
$ws.Range("AD6").Value = 5.05
$ws.Range("AE6").Value = 17092
$ws.Range("AF6").Value = 0.32
$ws.Range("AG6").Value = 75
$ws.Range("AH6").Value = 1.36
$ws.Range("AI6").Value = 6.79
$ws.Range("AJ6").Value = 53543977

# Rows 7-9: clear all data columns (D:AJ), keep only index/period columns A-C
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
